$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 800.86206
$ws.Range("I15").Value = 800.86206
$ws.Range("K15").Value = 2402.58618
$ws.Range("M15").Value = -2233.58618

# ALC row 40
$ws.Range("H40").Value = 1800
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1700
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1700
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2050

# ALC row 43
$ws.Range("H43").Value = 1273.7059
$ws.Range("I43").Value = 802
$ws.Range("J43").Value = 1693
$ws.Range("K43").Value = 802
$ws.Range("L43").Value = 1693
$ws.Range("M43").Value = -733
$ws.Range("N43").Value = -1831

# ALC row 125
$ws.Range("H125").Value = 3385.2856
$ws.Range("I125").Value = 2312.75
$ws.Range("J125").Value = 4815.3335
$ws.Range("K125").Value = 20814.75
$ws.Range("L125").Value = 43338.0015
$ws.Range("M125").Value = -18354.75
$ws.Range("N125").Value = -48258.0015

# ALC row 137
$ws.Range("H137").Value = 1490564
$ws.Range("I137").Value = 2646653.8
$ws.Range("K137").Value = 7939961.399999999
$ws.Range("M137").Value = -7937411.399999999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 4983.375
$ws.Range("I32").Value = 5009.4883
$ws.Range("J32").Value = 4929.905
$ws.Range("K32").Value = 5009.4883
$ws.Range("L32").Value = 4929.905
$ws.Range("M32").Value = -4722.4883
$ws.Range("N32").Value = -5503.905

# ARM row 76
$ws.Range("H76").Value = 30115.2
$ws.Range("J76").Value = 30115.2
$ws.Range("L76").Value = 30115.2
$ws.Range("N76").Value = -30791.2

# ARM row 79
$ws.Range("H79").Value = 30115.2
$ws.Range("J79").Value = 30115.2
$ws.Range("L79").Value = 30115.2
$ws.Range("N79").Value = -32455.2

# ARM row 102
$ws.Range("H102").Value = 3401.5
$ws.Range("J102").Value = 2999.5
$ws.Range("L102").Value = 2999.5
$ws.Range("N102").Value = -6243.5

# ARM row 135
$ws.Range("H135").Value = 40107
$ws.Range("J135").Value = 40107
$ws.Range("L135").Value = 40107
$ws.Range("N135").Value = -50247

# ARM row 137
$ws.Range("H137").Value = 43936.668
$ws.Range("J137").Value = 43936.668
$ws.Range("L137").Value = 43936.668
$ws.Range("N137").Value = -54136.668

$ws = $wb.Worksheets.Item("BSM")
# BSM row 59
$ws.Range("H59").Value = 42000
$ws.Range("J59").Value = 42000
$ws.Range("L59").Value = 42000
$ws.Range("N59").Value = -43694

# BSM row 80
$ws.Range("H80").Value = 200.34375
$ws.Range("I80").Value = 182
$ws.Range("K80").Value = 182
$ws.Range("M80").Value = 816

# BSM row 83
$ws.Range("H83").Value = 200.34375
$ws.Range("I83").Value = 182
$ws.Range("K83").Value = 910
$ws.Range("M83").Value = 4082

# BSM row 94
$ws.Range("H94").Value = 1182.5
$ws.Range("I94").Value = 853.125
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 853.125
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -402.125
$ws.Range("N94").Value = -3402

# BSM row 137
$ws.Range("H137").Value = 45350
$ws.Range("J137").Value = 45350
$ws.Range("L137").Value = 45350
$ws.Range("N137").Value = -55550

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58
$ws.Range("H58").Value = 3023.2666
$ws.Range("I58").Value = 1803.92
$ws.Range("K58").Value = 1803.92
$ws.Range("M58").Value = -1600.92

# CRP row 94
$ws.Range("H94").Value = 928.8889
$ws.Range("I94").Value = 633.0714
$ws.Range("J94").Value = 1247.4615
$ws.Range("K94").Value = 633.0714
$ws.Range("L94").Value = 1247.4615
$ws.Range("M94").Value = -182.0714
$ws.Range("N94").Value = -2149.4615

# CRP row 105
$ws.Range("H105").Value = 2685.2307
$ws.Range("I105").Value = 2802.25
$ws.Range("J105").Value = 2498
$ws.Range("K105").Value = 2802.25
$ws.Range("L105").Value = 2498
$ws.Range("M105").Value = -1055.25
$ws.Range("N105").Value = -5992

# CRP row 134
$ws.Range("H134").Value = 1848.5238
$ws.Range("I134").Value = 930.5294
$ws.Range("K134").Value = 2791.5882
$ws.Range("M134").Value = -256.5882000000001

# CRP row 136
$ws.Range("H136").Value = 3023.2666
$ws.Range("I136").Value = 1803.92
$ws.Range("K136").Value = 5411.76
$ws.Range("M136").Value = -2861.76

$ws = $wb.Worksheets.Item("CUL")
# CUL row 107
$ws.Range("H107").Value = 83770.914
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 167141.83
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 501425.49
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -505265.49

# CUL row 113
$ws.Range("H113").Value = 5435445
$ws.Range("I113").Value = 696.0833
$ws.Range("J113").Value = 11364262
$ws.Range("K113").Value = 2088.2499
$ws.Range("L113").Value = 34092786
$ws.Range("M113").Value = 81.7501000000002
$ws.Range("N113").Value = -34097126

# CUL row 122
$ws.Range("H122").Value = 2537.1633
$ws.Range("J122").Value = 3361.121
$ws.Range("L122").Value = 30250.089
$ws.Range("N122").Value = -35150.089

# CUL row 132
$ws.Range("H132").Value = 2204.6316
$ws.Range("I132").Value = 977.8570999999999
$ws.Range("K132").Value = 8800.713899999999
$ws.Range("M132").Value = -6270.713899999999

$ws = $wb.Worksheets.Item("GSM")
# GSM row 12
$ws.Range("H12").Value = 19870.6
$ws.Range("I12").Value = 5001.5
$ws.Range("J12").Value = 29783.334
$ws.Range("K12").Value = 5001.5
$ws.Range("L12").Value = 29783.334
$ws.Range("M12").Value = -4861.5
$ws.Range("N12").Value = -30063.334

# GSM row 137
$ws.Range("H137").Value = 34182.5
$ws.Range("J137").Value = 34182.5
$ws.Range("L137").Value = 34182.5
$ws.Range("N137").Value = -44382.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# LTW row 111
$ws.Range("H111").Value = 73695
$ws.Range("J111").Value = 73695
$ws.Range("L111").Value = 73695
$ws.Range("N111").Value = -81875

$ws = $wb.Worksheets.Item("WVR")
# WVR row 12
$ws.Range("H12").Value = 7250
$ws.Range("I12").Value = 6500
$ws.Range("J12").Value = 8000
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = -6358
$ws.Range("N12").Value = -8284

# WVR row 57
$ws.Range("H57").Value = 21300
$ws.Range("J57").Value = 21300
$ws.Range("L57").Value = 21300
$ws.Range("N57").Value = -22808

# WVR row 96
$ws.Range("H96").Value = 202100800
$ws.Range("I96").Value = 252625500
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 252625500
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -252624127
$ws.Range("N96").Value = -4746

# WVR row 126
$ws.Range("H126").Value = 2472.8572
$ws.Range("I126").Value = 1408.3334
$ws.Range("J126").Value = 4389
$ws.Range("K126").Value = 4225.0002
$ws.Range("L126").Value = 13167
$ws.Range("M126").Value = -1755.0002
$ws.Range("N126").Value = -18107
